# Apply the "Winner determination" summary-rows edit:
#  - add an average row for column J (|S*|/n) under the data block
#  - add four labelled summary rows (Average/Worst of SW and SC ratios)
#  - bold/size formatting on the new summary rows
#  - widen column A, set the printed page to A4 portrait, and leave the
#    selection on B18 (matches the state the workbook was saved in)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new average-of-column-J cell -----------------------------------------
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# --- four labelled summary rows --------------------------------------------
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(N2:N11)"

$summary = $ws.Range("A14:B17")
$summary.Font.Bold = $true
$summary.Font.Size = 12
$summary.VerticalAlignment = -4108

# --- cosmetic / layout tweaks matching the saved workbook ------------------
$ws.Columns.Item(1).ColumnWidth = 23.6

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B18").Select()
